$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 40
$ws.Range("D2").Value = 44238
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/caja 15 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 15

# Row 4 <- original row 32
$ws.Range("D4").Value = 44592
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 54
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("Q4").Value = '$/caja 15 kilos empedrada'
$ws.Range("R4").Value = 'Provincia de Limarí'
$ws.Range("S4").Value = 1333
$ws.Range("T4").Value = 15

# Row 5 <- original row 9
$ws.Range("D5").Value = 44271
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/caja 14 kilos granel'
$ws.Range("R5").Value = 'Provincia del Elquí'
$ws.Range("S5").Value = 857
$ws.Range("T5").Value = 14

# Row 6 <- original row 2
$ws.Range("D6").Value = 44270
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 85
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 14 kilos granel'
$ws.Range("R6").Value = 'Provincia del Elquí'
$ws.Range("S6").Value = 857
$ws.Range("T6").Value = 14

# Row 7 <- original row 31
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 48
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = '$/caja 14 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 14

# Row 8 <- original row 21
$ws.Range("D8").Value = 44313
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 36
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = '$/caja 14 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 14

# Row 9 <- original row 8
$ws.Range("D9").Value = 44614
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 54
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = '$/caja 14 kilos granel'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 14

# Row 10 <- original row 29
$ws.Range("D10").Value = 45044
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("Q10").Value = '$/caja 14 kilos granel'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 1143
$ws.Range("T10").Value = 14

# Row 11 <- original row 30
$ws.Range("D11").Value = 45044
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("Q11").Value = '$/caja 14 kilos granel'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 14

# Row 12 <- original row 41
$ws.Range("D12").Value = 44245
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 15

# Row 13 <- original row 18
$ws.Range("D13").Value = 44259
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = '$/caja 15 kilos empedrada'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 800
$ws.Range("T13").Value = 15

# Row 14 <- original row 4
$ws.Range("D14").Value = 45001
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("Q14").Value = '$/caja 14 kilos empedrada'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1143
$ws.Range("T14").Value = 14

# Row 15 <- original row 6
$ws.Range("D15").Value = 45042
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("Q15").Value = '$/caja 14 kilos granel'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 1214
$ws.Range("T15").Value = 14

# Row 16 <- original row 7
$ws.Range("D16").Value = 45042
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = '$/caja 14 kilos granel'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 14

# Row 17 <- original row 5
$ws.Range("D17").Value = 45015
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 56
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("Q17").Value = '$/caja 14 kilos empedrada'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 1071
$ws.Range("T17").Value = 14

# Row 18 <- original row 16
$ws.Range("D18").Value = 44260
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 56
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 13000
$ws.Range("Q18").Value = '$/caja 14 kilos empedrada'
$ws.Range("R18").Value = 'Provincia del Elquí'
$ws.Range("S18").Value = 929
$ws.Range("T18").Value = 14

# Row 19 <- original row 22
$ws.Range("D19").Value = 45043
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 17000
$ws.Range("P19").Value = 17000
$ws.Range("Q19").Value = '$/caja 14 kilos granel'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1214
$ws.Range("T19").Value = 14

# Row 20 <- original row 23
$ws.Range("D20").Value = 45043
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 67
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("Q20").Value = '$/caja 14 kilos granel'
$ws.Range("R20").Value = 'Provincia de Limarí'
$ws.Range("S20").Value = 1000
$ws.Range("T20").Value = 14

# Row 21 <- original row 14
$ws.Range("D21").Value = 45054
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 54
$ws.Range("N21").Value = 16000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 16000
$ws.Range("Q21").Value = '$/caja 14 kilos empedrada'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 1143
$ws.Range("T21").Value = 14

# Row 22 <- original row 15
$ws.Range("D22").Value = 45054
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("Q22").Value = '$/caja 14 kilos empedrada'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 14

# Row 23 <- original row 27
$ws.Range("D23").Value = 44278
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("Q23").Value = '$/caja 14 kilos empedrada'
$ws.Range("R23").Value = 'Provincia del Elquí'
$ws.Range("S23").Value = 929
$ws.Range("T23").Value = 14

# Row 24 <- original row 19
$ws.Range("D24").Value = 44315
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 65
$ws.Range("N24").Value = 14000
$ws.Range("O24").Value = 14000
$ws.Range("P24").Value = 14000
$ws.Range("Q24").Value = '$/caja 14 kilos granel'
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 1000
$ws.Range("T24").Value = 14

# Row 25 <- original row 28
$ws.Range("D25").Value = 44314
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 56
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("Q25").Value = '$/caja 14 kilos granel'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 14

# Row 26 <- original row 36
$ws.Range("D26").Value = 45040
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 65
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 17000
$ws.Range("P26").Value = 17000
$ws.Range("Q26").Value = '$/caja 14 kilos granel'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 1214
$ws.Range("T26").Value = 14

# Row 27 <- original row 37
$ws.Range("D27").Value = 45040
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 60
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("Q27").Value = '$/caja 14 kilos granel'
$ws.Range("R27").Value = 'Provincia de Limarí'
$ws.Range("S27").Value = 1000
$ws.Range("T27").Value = 14

# Row 28 <- original row 11
$ws.Range("D28").Value = 44323
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 14000
$ws.Range("O28").Value = 14000
$ws.Range("P28").Value = 14000
$ws.Range("Q28").Value = '$/caja 14 kilos granel'
$ws.Range("R28").Value = 'Provincia de Limarí'
$ws.Range("S28").Value = 1000
$ws.Range("T28").Value = 14

# Row 29 <- original row 39
$ws.Range("D29").Value = 45006
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 40
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 16000
$ws.Range("Q29").Value = '$/caja 14 kilos empedrada'
$ws.Range("R29").Value = 'Provincia de Limarí'
$ws.Range("S29").Value = 1143
$ws.Range("T29").Value = 14

# Row 30 <- original row 13
$ws.Range("D30").Value = 44312
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 68
$ws.Range("N30").Value = 14000
$ws.Range("O30").Value = 14000
$ws.Range("P30").Value = 14000
$ws.Range("Q30").Value = '$/caja 14 kilos granel'
$ws.Range("R30").Value = 'Provincia de Limarí'
$ws.Range("S30").Value = 1000
$ws.Range("T30").Value = 14

# Row 31 <- original row 10
$ws.Range("D31").Value = 44627
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 56
$ws.Range("N31").Value = 17000
$ws.Range("O31").Value = 17000
$ws.Range("P31").Value = 17000
$ws.Range("Q31").Value = '$/caja 14 kilos empedrada'
$ws.Range("R31").Value = 'Provincia de Limarí'
$ws.Range("S31").Value = 1214
$ws.Range("T31").Value = 14

# Row 32 <- original row 20
$ws.Range("D32").Value = 44320
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 45
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 14000
$ws.Range("Q32").Value = '$/caja 14 kilos granel'
$ws.Range("R32").Value = 'Provincia de Limarí'
$ws.Range("S32").Value = 1000
$ws.Range("T32").Value = 14

# Row 33 <- original row 26
$ws.Range("D33").Value = 44322
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 14000
$ws.Range("O33").Value = 14000
$ws.Range("P33").Value = 14000
$ws.Range("Q33").Value = '$/caja 14 kilos granel'
$ws.Range("R33").Value = 'Provincia de Limarí'
$ws.Range("S33").Value = 1000
$ws.Range("T33").Value = 14

# Row 34 <- original row 33
$ws.Range("D34").Value = 44630
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 75
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 15000
$ws.Range("Q34").Value = '$/caja 14 kilos empedrada'
$ws.Range("R34").Value = 'Provincia de Limarí'
$ws.Range("S34").Value = 1071
$ws.Range("T34").Value = 14

# Row 35 <- original row 24
$ws.Range("D35").Value = 45050
$ws.Range("L35").Value = 'Especial'
$ws.Range("M35").Value = 56
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 14000
$ws.Range("P35").Value = 14000
$ws.Range("Q35").Value = '$/caja 14 kilos granel'
$ws.Range("R35").Value = 'Provincia de Limarí'
$ws.Range("S35").Value = 1000
$ws.Range("T35").Value = 14

# Row 36 <- original row 25
$ws.Range("D36").Value = 45050
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 12000
$ws.Range("P36").Value = 12000
$ws.Range("Q36").Value = '$/caja 14 kilos granel'
$ws.Range("R36").Value = 'Provincia de Limarí'
$ws.Range("S36").Value = 857
$ws.Range("T36").Value = 14

# Row 37 <- original row 35
$ws.Range("D37").Value = 44616
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 70
$ws.Range("N37").Value = 14000
$ws.Range("O37").Value = 14000
$ws.Range("P37").Value = 14000
$ws.Range("Q37").Value = '$/caja 14 kilos empedrada'
$ws.Range("R37").Value = 'Provincia de Limarí'
$ws.Range("S37").Value = 1000
$ws.Range("T37").Value = 14

# Row 38 <- original row 34
$ws.Range("D38").Value = 45014
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 60
$ws.Range("N38").Value = 15000
$ws.Range("O38").Value = 15000
$ws.Range("P38").Value = 15000
$ws.Range("Q38").Value = '$/caja 14 kilos empedrada'
$ws.Range("R38").Value = 'Provincia de Limarí'
$ws.Range("S38").Value = 1071
$ws.Range("T38").Value = 14

# Row 39 <- original row 12
$ws.Range("D39").Value = 44239
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 70
$ws.Range("N39").Value = 15000
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 15000
$ws.Range("Q39").Value = '$/caja 15 kilos granel'
$ws.Range("R39").Value = 'Provincia de Limarí'
$ws.Range("S39").Value = 1000
$ws.Range("T39").Value = 15

# Row 40 <- original row 38
$ws.Range("D40").Value = 44242
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 45
$ws.Range("N40").Value = 12000
$ws.Range("O40").Value = 12000
$ws.Range("P40").Value = 12000
$ws.Range("Q40").Value = '$/caja 15 kilos granel'
$ws.Range("R40").Value = 'Provincia de Limarí'
$ws.Range("S40").Value = 800
$ws.Range("T40").Value = 15

# Row 41 <- original row 17
$ws.Range("D41").Value = 44588
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 85
$ws.Range("N41").Value = 19000
$ws.Range("O41").Value = 20000
$ws.Range("P41").Value = 19529
$ws.Range("Q41").Value = '$/caja 14 kilos granel'
$ws.Range("R41").Value = 'Provincia de Limarí'
$ws.Range("S41").Value = 1395
$ws.Range("T41").Value = 14

